# Generate Report for Handoff
#
# Updates the localization-status report so the row for
# f6dc4ae2-9dd8-455d-8905-b66f00492b6a.md reflects that the handback file
# has been regenerated and is ready for handoff (status/date/error-detail
# refresh on the "Overview", "zh-cn" and "de-de" sheets).

$wb = $excel.ActiveWorkbook

$errorDetail = 'The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/18fdb801e619b1c2eef5d93a59f4c63d964f8ee3/e2e/f6dc4ae2-9dd8-455d-8905-b66f00492b6a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/54434be78519ec2f1b42d2dc71e4875769396e2b/e2e/f6dc4ae2-9dd8-455d-8905-b66f00492b6a.md.'

# ---- Overview sheet: row 3 is the f6dc4ae2-...md file ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-27 04:46:47"

# ---- zh-cn sheet: row 3 is the f6dc4ae2-...md file ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("H3").Value = "2016-08-27 04:46:43"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = 39.1666666666667

# ---- de-de sheet: row 3 is the f6dc4ae2-...md file ----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("H3").Value = "2016-08-27 04:46:47"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 39.1666666666667
